$d = $word.ActiveDocument

# Locate the paragraph that starts with Paul's comment about Figure 6.
# We search using a substring that avoids "smart quote" characters (the
# source text contains a curly apostrophe in "don't") so the Find engine
# matches reliably.
$rng = $d.Content
$found = $rng.Find.Execute("Figure 6 needs a plain language title or label indicating", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph for highlighting"
}

# Walk the document paragraphs and apply a green highlight (matching
# wdBrightGreen / OOXML w:highlight w:val="green") to the whole paragraph
# -- both the run text and the paragraph mark -- for the paragraph that
# contains the located text.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Figure 6 needs a plain language title*") {
        $p.Range.Font.HighlightColorIndex = 4
    }
}
